$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Append new row 11 with the same structure as the preceding rows ---
$ws.Cells.Item(11, 1).Value = "Demo inplannen"
$ws.Cells.Item(11, 2).Value = "klantenservice@testbedrijf123.nl"
$ws.Cells.Item(11, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Cells.Item(11, 4).Value = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item(11, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Cells.Item(11, 6).Value = "2025-08-13 21:28:22"
$ws.Cells.Item(11, 7).Value = "Nee"
$ws.Cells.Item(11, 8).Value = "Ja"
$ws.Cells.Item(11, 9).Value = "Nee"
$ws.Cells.Item(11, 10).Value = "Nee"

# --- Extend the conditional-formatting ranges so the new row is covered ---
$ws.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D11"))
$ws.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G11"))
$ws.Range("H2:H10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H11"))
$ws.Range("I2:I10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I11"))
$ws.Range("J2:J10").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J11"))

# --- Update the Dashboard summary count for this category ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 10
